# Remove the direct "2B2B2B" font-color override from the paragraph mark
# (w:pPr/w:rPr) and from the single run (w:r/w:rPr) of the
# "quarterly change ($) (2.5 points)" bullet, matching the commit's XML
# diff. Word's object model has no notion of literally deleting the
# <w:color> element, so we restore it to the automatic color (the
# standard "no manual color override" state) which is how Word itself
# represents a cleared font color.

$d = $word.ActiveDocument

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Text = "quarterly change ($) (2.5 points)"
$rng.Find.MatchWildcards = $false
$rng.Find.MatchCase = $true
$found = $rng.Find.Execute()

if ($found) {
    # Grab the whole paragraph (this includes the trailing paragraph
    # mark), so the color change lands on both the run's rPr and the
    # paragraph mark's rPr inside pPr.
    $para = $rng.Paragraphs(1)
    $prange = $para.Range
    $prange.Font.Color = -16777216
}

Write-Output "done"
